# Apply the edits described by the diff:
#  - Add a new row (Sita / 3500) to the "Seattle" worksheet (sheet3.xml)
#  - Make "Seattle" the active/selected sheet, with H14 selected
#  - The new string "Sita" is automatically appended to the shared
#    strings table when the cell value is written

$wb = $excel.ActiveWorkbook

$seattle = $wb.Worksheets.Item("Seattle")

# Switch to the Seattle sheet - this becomes the active tab,
# which also clears tabSelected on the previously active sheet (CA)
$seattle.Activate()

# Add the new data row under the existing data (row 3)
$seattle.Range("A3").Value = "Sita"
$seattle.Range("B3").Value = 3500

# Update the selection on the Seattle sheet to H14
$seattle.Range("H14").Select()
